$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 14 (last row removed entirely)
$ws.Rows("14").Delete()

# Step 2a: style fixes where donor cell itself will also change style later
#          (must copy BEFORE the donor's own style is overwritten)
$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("D11").PasteSpecial(-4122)

# Step 2b: remaining style fixes using stable donor cells
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: update cell text/values to match the after-state
$ws.Range("B2").Value = "MOV1_2`n(D0,R0,R0)"
$ws.Range("C2").Value = "OX(U0,L0,D0)"
$ws.Range("D2").Value = "OX(U0,R0,D0)"
$ws.Range("B3").Value = "ADD_DATA`n(D1,U0)`nDATA(0x1)"
$ws.Range("C3").Value = "SWITCH_TAG`n(U0,U1,D0,D1)`nDATA(0x0100)"
$ws.Range("D3").Value = "SWITCH_TAG`n(U0,U1,D0,D1)`nDATA(0x0100)"
$ws.Range("E3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "SYNC`n(0x01,U0,L0,R0,D0)`nDATA(0x0)"
$ws.Range("D4").Value = "SYNC`n(0x01,U0,R0,L0,D0)`nDATA(0x0)"
$ws.Range("A5").Value = "DISCARD`n(R0,R0,R0,R0)"
$ws.Range("B5").Value = "SWITCH_PRED`n(D1,D0,U1,L0)"
$ws.Range("C5").Value = "SWITCH_TAG`n(U1,U0,D0,D1)`nDATA(0x0200)"
$ws.Range("D5").Value = "SWITCH_TAG`n(U1,U0,D0,D1)`nDATA(0x0200)"
$ws.Range("E5").ClearContents()
$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "MOV1_2`n(D1,U0,U0)"
$ws.Range("C6").Value = "SYNC`n(0x01,U0,L0,R0,D0)`nDATA(0x0)"
$ws.Range("D6").Value = "SYNC`n(0x01,U0,R0,L0,D0)`nDATA(0x0)"
$ws.Range("B7").Value = "MOV1_2`n(D1,U1,U1)"
$ws.Range("C7").Value = "SWITCH_TAG`n(U1,U0,D0,D1)`nDATA(0x0300)"
$ws.Range("D7").Value = "SWITCH_TAG`n(U1,U0,D0,D1)`nDATA(0x0300)"
$ws.Range("E7").Value = "DISCARD`n(D0,D0,D0,D0)"
$ws.Range("C8").Value = "SYNC`n(0x01,U0,L0,R0,D0)`nDATA(0x0)"
$ws.Range("D8").Value = "SYNC`n(0x01,U0,R0,L0,D0)`nDATA(0x0)"
$ws.Range("E8").Value = "SWITCH_PRED`n(D1,D0,U1,U0)"
$ws.Range("C9").Value = "SWITCH_TAG`n(U1,U0,D0,D1)`nDATA(0x0400)"
$ws.Range("D9").Value = "SWITCH_TAG`n(U1,U0,D0,D1)`nDATA(0x0400)"
$ws.Range("E9").Value = "MOV1_2`n(D1,U0,U0)"
$ws.Range("C10").Value = "SYNC`n(0x01,U0,L0,R0,D0)`nDATA(0x0)"
$ws.Range("D10").Value = "SYNC`n(0x01,U0,R0,L0,D0)`nDATA(0x0)"
$ws.Range("E10").Value = "AND`n(D0,D1,U1)"
$ws.Range("B11").Value = "MOV1_2`n(R0,U1,U1)"
$ws.Range("C11").Value = "MOV2_1`n(U1,U0,L0)"
$ws.Range("D11").Value = "MOV2_1`n(U1,U0,D0)"
$ws.Range("E11").Value = "ADD_DATA`n(D0,U0)`nMOV1_2`n(D1,U1,U1)`nDATA(0xff)"
$ws.Range("B12").Value = "MOV1_2`n(D0,U1,U1)"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "MOV1_2`n(U0,D0,R0)"
$ws.Range("E12").Value = "MOV1_2`n(L0,U0,U1)"
$ws.Range("B13").Value = "MOV1_2`n(R0,U0,U0)"
$ws.Range("C13").Value = "MOV1_2`n(R1,L0,L0)"
$ws.Range("D13").Value = "NE_DATA`n(U0,R0)`nDATA(0x0)"
$ws.Range("E13").Value = "MOV1_2`n(L0,U1,L1)"

# Step 4: selection to match the after-state sheetView
$ws.Range("G2").Select()
